# Weekly refresh of the "Hortaliza, Vega Monumental Concepción - Haba" sheet.
# The commit reshuffles each data row's Fecha / Volumen / Precio mínimo /
# Precio máximo / Precio promedio ponderado / Origen / Precio $/Kg
# (columns D, J, K, L, M, O, P) to a new weekly snapshot, while the rest of
# the row (Mercado, Región, Categoría, Variedad, Calidad, Unidad de
# comercialización, Kg o Unidades, Clasificación) stays put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row, Fecha(D), Volumen(J), PrecioMin(K), PrecioMax(L), PrecioProm(M), Origen(O), PrecioKg(P)
$rowsData = @(
    @(2, 44351, 100, 15000, 16000, 15500, 'Región Metropolitana', 620),
    @(3, 44497, 150, 6000, 6500, 6333, 'Región Metropolitana', 253),
    @(4, 44537, 160, 8500, 9000, 8719, 'Región del Maule', 349),
    @(5, 44188, 100, 18000, 20000, 19000, 'Región Metropolitana', 760),
    @(6, 44540, 140, 11000, 12000, 11429, 'Región del Maule', 457),
    @(7, 44523, 100, 9000, 10000, 9500, 'Región Metropolitana', 380),
    @(8, 44505, 180, 6000, 6500, 6222, 'Región del Maule', 249),
    @(9, 44498, 220, 7000, 7500, 7273, 'Región Metropolitana', 291),
    @(10, 44160, 100, 9000, 10000, 9500, 'Región Metropolitana', 380),
    @(11, 44335, 100, 18000, 20000, 19000, 'Provincia de Limarí', 760),
    @(12, 44316, 100, 16000, 18000, 17000, 'Región Metropolitana', 680),
    @(13, 44526, 100, 7500, 8000, 7750, 'Región Metropolitana', 310),
    @(14, 44467, 100, 8000, 9000, 8500, 'Región Metropolitana', 340),
    @(15, 44517, 130, 6000, 6500, 6269, 'Región Metropolitana', 251),
    @(16, 44483, 350, 5500, 6000, 5714, 'Región Metropolitana', 229),
    @(17, 44476, 100, 7000, 7500, 7250, 'Región Metropolitana', 290),
    @(18, 44509, 100, 6500, 7000, 6750, 'Región Metropolitana', 270),
    @(19, 44482, 430, 8000, 8500, 8267, "Región de O'Higgins", 331),
    @(20, 44461, 100, 13000, 14000, 13500, 'Provincia del Elquí', 540),
    @(21, 44545, 140, 14000, 15000, 14429, 'Provincia de Chacabuco', 577),
    @(22, 44454, 100, 13000, 14000, 13500, 'Provincia del Elquí', 540),
    @(23, 44162, 100, 7500, 8000, 7750, 'Región Metropolitana', 310),
    @(24, 44533, 180, 8000, 8500, 8222, 'Región del Maule', 329),
    @(25, 44503, 250, 9000, 10000, 9400, 'Provincia de Melipilla', 376),
    @(26, 44384, 100, 12000, 13000, 12500, 'Región de Coquimbo', 500)
)

foreach ($entry in $rowsData) {
    $r = $entry[0]
    $ws.Range("D$r").Value2 = $entry[1]
    $ws.Range("J$r").Value2 = $entry[2]
    $ws.Range("K$r").Value2 = $entry[3]
    $ws.Range("L$r").Value2 = $entry[4]
    $ws.Range("M$r").Value2 = $entry[5]
    $ws.Range("O$r").Value2 = $entry[6]
    $ws.Range("P$r").Value2 = $entry[7]
}
